$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at row 27 for the new "31/12/2023" data block
$ws.Range("A27:D31").EntireRow.Insert()

$ws.Cells.Item(27,1).Value = "Brasil"
$ws.Cells.Item(27,2).Value = "Geladeira"
$ws.Cells.Item(27,3).Value = "31/12/2023"
$ws.Cells.Item(27,4).Value = 98.2
$ws.Cells.Item(28,1).Value = "Brasil"
$ws.Cells.Item(28,2).Value = "Máquina de lavar roupa"
$ws.Cells.Item(28,3).Value = "31/12/2023"
$ws.Cells.Item(28,4).Value = 69.4
$ws.Cells.Item(29,1).Value = "Brasil"
$ws.Cells.Item(29,2).Value = "Carro"
$ws.Cells.Item(29,3).Value = "31/12/2023"
$ws.Cells.Item(29,4).Value = 48.1
$ws.Cells.Item(30,1).Value = "Brasil"
$ws.Cells.Item(30,2).Value = "Motocicleta"
$ws.Cells.Item(30,3).Value = "31/12/2023"
$ws.Cells.Item(30,4).Value = 24.6
$ws.Cells.Item(31,1).Value = "Brasil"
$ws.Cells.Item(31,2).Value = "Carro e motocicleta"
$ws.Cells.Item(31,3).Value = "31/12/2023"
$ws.Cells.Item(31,4).Value = 12.6

# Insert 5 new rows at row 57 for the new "31/12/2023" data block
$ws.Range("A57:D61").EntireRow.Insert()

$ws.Cells.Item(57,1).Value = "Nordeste"
$ws.Cells.Item(57,2).Value = "Geladeira"
$ws.Cells.Item(57,3).Value = "31/12/2023"
$ws.Cells.Item(57,4).Value = 97
$ws.Cells.Item(58,1).Value = "Nordeste"
$ws.Cells.Item(58,2).Value = "Máquina de lavar roupa"
$ws.Cells.Item(58,3).Value = "31/12/2023"
$ws.Cells.Item(58,4).Value = 39
$ws.Cells.Item(59,1).Value = "Nordeste"
$ws.Cells.Item(59,2).Value = "Carro"
$ws.Cells.Item(59,3).Value = "31/12/2023"
$ws.Cells.Item(59,4).Value = 27.6
$ws.Cells.Item(60,1).Value = "Nordeste"
$ws.Cells.Item(60,2).Value = "Motocicleta"
$ws.Cells.Item(60,3).Value = "31/12/2023"
$ws.Cells.Item(60,4).Value = 32.6
$ws.Cells.Item(61,1).Value = "Nordeste"
$ws.Cells.Item(61,2).Value = "Carro e motocicleta"
$ws.Cells.Item(61,3).Value = "31/12/2023"
$ws.Cells.Item(61,4).Value = 9.7

# Insert 5 new rows at row 87 for the new "31/12/2023" data block
$ws.Range("A87:D91").EntireRow.Insert()

$ws.Cells.Item(87,1).Value = "Sergipe"
$ws.Cells.Item(87,2).Value = "Geladeira"
$ws.Cells.Item(87,3).Value = "31/12/2023"
$ws.Cells.Item(87,4).Value = 97.4
$ws.Cells.Item(88,1).Value = "Sergipe"
$ws.Cells.Item(88,2).Value = "Máquina de lavar roupa"
$ws.Cells.Item(88,3).Value = "31/12/2023"
$ws.Cells.Item(88,4).Value = 41
$ws.Cells.Item(89,1).Value = "Sergipe"
$ws.Cells.Item(89,2).Value = "Carro"
$ws.Cells.Item(89,3).Value = "31/12/2023"
$ws.Cells.Item(89,4).Value = 27.9
$ws.Cells.Item(90,1).Value = "Sergipe"
$ws.Cells.Item(90,2).Value = "Motocicleta"
$ws.Cells.Item(90,3).Value = "31/12/2023"
$ws.Cells.Item(90,4).Value = 28.4
$ws.Cells.Item(91,1).Value = "Sergipe"
$ws.Cells.Item(91,2).Value = "Carro e motocicleta"
$ws.Cells.Item(91,3).Value = "31/12/2023"
$ws.Cells.Item(91,4).Value = 8.2
